$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("food", "食べ物|たべもの")
    ,@("drink", "飲み物|のみもの")
    ,@("things", "物|もの")
    ,@("shopping", "買い物|かいもの")
    ,@("animal", "動物|どうぶつ")
    ,@("commodity prices", "物価|ぶっか")
    ,@("bird", "鳥|とり")
    ,@("grilled chicken", "焼き鳥|やきとり")
    ,@("swan", "白鳥|はくちょう")
    ,@("cooking", "料理|りょうり")
    ,@("charge", "料金|りょうきん")
    ,@("tuition", "授業料|じゅぎょうりょう")
    ,@("salary", "給料|きゅうりょう")
    ,@("cooking", "料理|りょうり")
    ,@("reason", "理由|りゆう")
    ,@("geography", "地理|ちり")
    ,@("impossible", "無理な|むりな")
    ,@("especially", "特に|とくに")
    ,@("special", "特別な|とくべつな")
    ,@("characteristic", "特徴|とくちょう")
    ,@("super express", "特急|とっきゅう")
    ,@("cheap", "安い|やすい")
    ,@("safe", "安全な|あんぜんな")
    ,@("relief", "安心|あんしん")
    ,@("anxious; worried", "不安な|ふあんな")
    ,@("rice; meal", "ご飯|ごはん")
    ,@("breakfast", "朝ご飯|あさごはん")
    ,@("lunch", "昼ご飯|ひるごはん")
    ,@("dinner", "晩ご飯|ばんごはん")
    ,@("meat", "肉|にく")
    ,@("beef", "牛肉|ぎゅうにく")
    ,@("pork", "豚肉|ぶたにく")
    ,@("meat shop", "肉屋|にくや")
    ,@("muscle", "筋肉|きんにく")
    ,@("bad", "悪い|わるい")
    ,@("to feel sick", "気分が悪い|きぶんがわるい")
    ,@("the worst", "最悪|さいあく")
    ,@("devil", "悪魔|あくま")
    ,@("body", "体|からだ")
    ,@("body weight", "体重|たいじゅう")
    ,@("gymnastics; physical exercises", "体操|たいそう")
    ,@("the same", "同じ|おなじ")
    ,@("coworker", "同僚|どうりょう")
    ,@("classmate", "同級生|どうきゅうせい")
    ,@("same time", "同時|どうじ")
    ,@("to arrive", "着く|つく")
    ,@("to wear", "着る|きる")
    ,@("kimono", "着物|きもの")
    ,@("swimwear", "水着|みずぎ")
    ,@("arriving at Osaka", "大阪着|おおさかちゃく")
    ,@("airport", "空港|くうこう")
    ,@("air", "空気|くうき")
    ,@("sky", "空|そら")
    ,@("to be vacant", "空く|あく")
    ,@("karate", "空手|からて")
    ,@("airport", "空港|くうこう")
    ,@("Kobe Port", "神戸港|こうべこう")
    ,@("port", "港|みなと")
    ,@("Hong Kong", "香港|ほんこん")
    ,@("noon; daytime", "昼|ひる")
    ,@("lunch", "昼ご飯|ひるごはん")
    ,@("nap", "昼寝|ひるね")
    ,@("lunch break", "昼休み|ひるやすみ")
    ,@("lunch (formal)", "昼食|ちゅうしょく")
    ,@("sea", "海|うみ")
    ,@("the Japan Sea", "日本海|にほんかい")
    ,@("overseas", "海外|かいがい")
    ,@("coast", "海岸|かいがん")
    ,@("Hokkaido", "北海道|ほっかいどう")
)

$startRow = 79
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
